$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the old row 2 (which held "0004 : Michele Hall"),
# pushing the existing client rows + Total row down from rows 2-4 to rows 5-7.
$ws.Range("A2:A4").EntireRow.Insert()

# New client rows inserted at rows 2-4.
$ws.Range("A2").Value = "0001 : Judith Lynch"
$ws.Range("B2").Value = 384.79

$ws.Range("A3").Value = "0002 : Lurline Odriscoll"
$ws.Range("B3").Value = 188.24

$ws.Range("A4").Value = "0003 : June Scala"
$ws.Range("B4").Value = 205.52

# Update the Total row (shifted down to row 7) to reflect the new sum.
$ws.Range("B7").Value = 982.45
